$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Insert two new columns before column D, shifting old D:K to F:M
$ws.Range("D:E").Insert()

# Step 2: Copy number formatting from column F (old column D, now shifted) into new columns D and E
# so the new quarters keep the same date/number formatting as the rest of the table.
$ws.Range("F5:F102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$ws.Range("F5:F102").Copy()
$ws.Range("E5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Step 3: Populate the two new quarter columns (D = Dec-2018, E = Sep-2018) with their reported figures
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 114900
$ws.Range("E8").Value = 119100
$ws.Range("D9").Value = 37600
$ws.Range("E9").Value = 37900
$ws.Range("D10").Value = 77300
$ws.Range("E10").Value = 81200
$ws.Range("D12").Value = 3000
$ws.Range("E12").Value = 2900
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = "NA"
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 72300
$ws.Range("E17").Value = 90400
$ws.Range("D18").Value = 42600
$ws.Range("E18").Value = 28700
$ws.Range("D20").Value = 4400
$ws.Range("E20").Value = 5400
$ws.Range("D21").Value = 52100
$ws.Range("E21").Value = 40600
$ws.Range("D22").Value = 100
$ws.Range("E22").Value = 100
$ws.Range("D23").Value = 46900
$ws.Range("E23").Value = 34000
$ws.Range("D24").Value = 8200
$ws.Range("E24").Value = 3900
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 38700
$ws.Range("E26").Value = 30100
$ws.Range("D27").Value = 31800
$ws.Range("E27").Value = 25400
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 7500
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -4400
$ws.Range("E32").Value = -5400
$ws.Range("D33").Value = 31800
$ws.Range("E33").Value = 32900
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 31800
$ws.Range("E35").Value = 32900
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 338900
$ws.Range("E41").Value = 186400
$ws.Range("D42").Value = 613500
$ws.Range("E42").Value = 766800
$ws.Range("D43").Value = 125100
$ws.Range("E43").Value = 117800
$ws.Range("D44").Value = 243300
$ws.Range("E44").Value = 234100
$ws.Range("D45").Value = 36400
$ws.Range("E45").Value = 29500
$ws.Range("D46").Value = 1357200
$ws.Range("E46").Value = 1334600
$ws.Range("D47").Value = 66200
$ws.Range("E47").Value = 56000
$ws.Range("D48").Value = 210500
$ws.Range("E48").Value = 205100
$ws.Range("D49").Value = 366800
$ws.Range("E49").Value = 367200
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 9200
$ws.Range("E52").Value = 11400
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 2010000
$ws.Range("E54").Value = 1974400
$ws.Range("D57").Value = 111300
$ws.Range("E57").Value = 112400
$ws.Range("D58").Value = "NA"
$ws.Range("E58").Value = "NA"
$ws.Range("D59").Value = 11000
$ws.Range("E59").Value = 11800
$ws.Range("D60").Value = 122300
$ws.Range("E60").Value = 124300
$ws.Range("D61").Value = 0
$ws.Range("E61").Value = 0
$ws.Range("D62").Value = 42900
$ws.Range("E62").Value = 42800
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 287900
$ws.Range("E66").Value = 285500
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 634500
$ws.Range("E72").Value = 599600
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 1722000
$ws.Range("E76").Value = 1688800
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 31800
$ws.Range("E81").Value = 32900
$ws.Range("D83").Value = 5000
$ws.Range("E83").Value = 6500
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 32900
$ws.Range("E89").Value = 25500
$ws.Range("D91").Value = -5600
$ws.Range("E91").Value = -7700
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = 127300
$ws.Range("E94").Value = -517300
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -10000
$ws.Range("E100").Value = 580500
$ws.Range("D101").Value = 2300
$ws.Range("E101").Value = -5600
$ws.Range("D102").Value = 152500
$ws.Range("E102").Value = 83100

# Step 4: A handful of previously-reported quarters were revised/corrected; update those cells too
$ws.Range("F48").Value = 211600
$ws.Range("G48").Value = 222400
$ws.Range("H48").Value = 191700
$ws.Range("I48").Value = 179700
$ws.Range("J48").Value = 169600
$ws.Range("F49").Value = 389000
$ws.Range("G49").Value = 391400
$ws.Range("H49").Value = "NA"
$ws.Range("I49").Value = "NA"
$ws.Range("J49").Value = "NA"
$ws.Range("F91").Value = -7400
$ws.Range("G91").Value = -11000
$ws.Range("H91").Value = -9700
$ws.Range("J91").Value = -7000
